# Apply updated "dSF" (column F) values for rows 2-16 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = -9
    4  = -3
    5  = -7
    7  = -6
    9  = 1
    10 = -10
    11 = -2
    12 = -5
    13 = -9
    15 = 9
    16 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
